$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.273.86"
$ws.Range("E2").Value = "  +3.02%  "

$ws.Range("D3").Value = "1.902.62"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5142"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.264"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "1.903.88"
$ws.Range("E13").Value = "  +1.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.360"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.34%  "

$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.040"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").Value = "29.285.01"
$ws.Range("E23").Value = "  +2.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  -2.85%  "

$ws.Range("D26").Value = "2.119.99"
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.456"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1049"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.056"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.647"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06598"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2199"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.238"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6523"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.232"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6056"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.678"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.057"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.231"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.67%  "
